# Auto-generated edit script applying the Lich_Profits market-data refresh diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) on affected rows
# across all 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

# --- ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H43").Value = 2999.4
$ws.Range("I43").Value = 2749.25
$ws.Range("K43").Value = 2749.25
$ws.Range("M43").Value = -2680.25
$ws.Range("H62").Value = 6660.5884
$ws.Range("I62").Value = 6532.6924
$ws.Range("J62").Value = 7076.25
$ws.Range("K62").Value = 6532.6924
$ws.Range("L62").Value = 7076.25
$ws.Range("M62").Value = -5908.6924
$ws.Range("N62").Value = -8324.25
$ws.Range("H65").Value = 6660.5884
$ws.Range("I65").Value = 6532.6924
$ws.Range("J65").Value = 7076.25
$ws.Range("K65").Value = 32663.462
$ws.Range("L65").Value = 35381.25
$ws.Range("M65").Value = -29543.462
$ws.Range("N65").Value = -41621.25
$ws.Range("H80").Value = 535.4483
$ws.Range("J80").Value = 813.3333
$ws.Range("L80").Value = 2439.9999
$ws.Range("N80").Value = -4435.9999
$ws.Range("H83").Value = 535.4483
$ws.Range("J83").Value = 813.3333
$ws.Range("L83").Value = 7319.9997
$ws.Range("N83").Value = -17303.9997
$ws.Range("H88").Value = 2178.8
$ws.Range("I88").Value = 1848.75
$ws.Range("J88").Value = 3499
$ws.Range("K88").Value = 1848.75
$ws.Range("L88").Value = 3499
$ws.Range("M88").Value = -1442.75
$ws.Range("N88").Value = -4311
$ws.Range("H91").Value = 2178.8
$ws.Range("I91").Value = 1848.75
$ws.Range("J91").Value = 3499
$ws.Range("K91").Value = 1848.75
$ws.Range("L91").Value = 3499
$ws.Range("M91").Value = -444.75
$ws.Range("N91").Value = -6307
$ws.Range("H125").Value = 2374.5
$ws.Range("I125").Value = 2229.5715
$ws.Range("K125").Value = 20066.1435
$ws.Range("M125").Value = -17606.1435
$ws.Range("H135").Value = 1453.56
$ws.Range("I135").Value = 1585.35
$ws.Range("K135").Value = 14268.15
$ws.Range("M135").Value = -11733.15

# --- ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()
$ws.Range("H32").Value = 10532.141
$ws.Range("I32").Value = 9906.638000000001
$ws.Range("K32").Value = 9906.638000000001
$ws.Range("M32").Value = -9619.638000000001
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H88").Value = 1788.6
$ws.Range("I88").Value = 1810.5
$ws.Range("K88").Value = 1810.5
$ws.Range("M88").Value = -1404.5
$ws.Range("H91").Value = 1788.6
$ws.Range("I91").Value = 1810.5
$ws.Range("K91").Value = 1810.5
$ws.Range("M91").Value = -406.5

# --- BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H36").Value = 4897.5
$ws.Range("I36").Value = 1795
$ws.Range("K36").Value = 1795
$ws.Range("M36").Value = -1261
$ws.Range("H86").Value = 1429.0667
$ws.Range("I86").Value = 1140.0454
$ws.Range("J86").Value = 2223.875
$ws.Range("K86").Value = 1140.0454
$ws.Range("L86").Value = 2223.875
$ws.Range("M86").Value = -17.04539999999997
$ws.Range("N86").Value = -4469.875
$ws.Range("H89").Value = 1429.0667
$ws.Range("I89").Value = 1140.0454
$ws.Range("J89").Value = 2223.875
$ws.Range("K89").Value = 5700.227
$ws.Range("L89").Value = 11119.375
$ws.Range("M89").Value = -84.22699999999986
$ws.Range("N89").Value = -22351.375
$ws.Range("H99").Value = 7224.606
$ws.Range("I99").Value = 7852.5835
$ws.Range("K99").Value = 7852.5835
$ws.Range("M99").Value = -6354.5835
$ws.Range("H134").Value = 1701.15
$ws.Range("I134").Value = 1701.15
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5103.450000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2568.450000000001
$ws.Range("N134").ClearContents()

# --- CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H105").Value = 4650.838
$ws.Range("I105").Value = 1771.75
$ws.Range("K105").Value = 1771.75
$ws.Range("M105").Value = -24.75
$ws.Range("H120").Value = 59344.4
$ws.Range("J120").Value = 61955.5
$ws.Range("L120").Value = 61955.5
$ws.Range("N120").Value = -69213.5
$ws.Range("H121").Value = 57530
$ws.Range("J121").Value = 59912.5
$ws.Range("L121").Value = 59912.5
$ws.Range("N121").Value = -62532.5

# --- CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H132").Value = 5744.5713
$ws.Range("I132").Value = 6680.136
$ws.Range("J132").Value = 2314.1667
$ws.Range("K132").Value = 60121.224
$ws.Range("L132").Value = 20827.5003
$ws.Range("M132").Value = -57591.224
$ws.Range("N132").Value = -25887.5003

# --- GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H3").Value = 1857.8
$ws.Range("I3").Value = 2266.3333
$ws.Range("J3").Value = 1245
$ws.Range("K3").Value = 2266.3333
$ws.Range("L3").Value = 1245
$ws.Range("M3").Value = -2150.3333
$ws.Range("N3").Value = -1477
$ws.Range("H80").Value = 3421.182
$ws.Range("I80").Value = 3092.6667
$ws.Range("K80").Value = 3092.6667
$ws.Range("M80").Value = -2094.6667
$ws.Range("H83").Value = 3421.182
$ws.Range("I83").Value = 3092.6667
$ws.Range("K83").Value = 15463.3335
$ws.Range("M83").Value = -10471.3335

# --- LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 2572.8484
$ws.Range("I68").Value = 2543.875
$ws.Range("K68").Value = 2543.875
$ws.Range("M68").Value = -1794.875
$ws.Range("H71").Value = 2572.8484
$ws.Range("I71").Value = 2543.875
$ws.Range("K71").Value = 12719.375
$ws.Range("M71").Value = -8975.375
$ws.Range("H82").Value = 1366.5454
$ws.Range("I82").Value = 1139.3334
$ws.Range("J82").Value = 1639.2
$ws.Range("K82").Value = 1139.3334
$ws.Range("L82").Value = 1639.2
$ws.Range("M82").Value = -778.3334
$ws.Range("N82").Value = -2361.2
$ws.Range("H85").Value = 1366.5454
$ws.Range("I85").Value = 1139.3334
$ws.Range("J85").Value = 1639.2
$ws.Range("K85").Value = 1139.3334
$ws.Range("L85").Value = 1639.2
$ws.Range("M85").Value = 108.6666
$ws.Range("N85").Value = -4135.2

# --- WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H20").Value = 100000000
$ws.Range("I20").Value = 100000000
$ws.Range("K20").Value = 100000000
$ws.Range("M20").Value = -99999760
$ws.Range("H135").Value = 89800
$ws.Range("J135").Value = 89800
$ws.Range("L135").Value = 89800
$ws.Range("N135").Value = -99940

